# Project "Sample Project" is saved. Update the rule name in cell B11 on
# the "Rules" sheet from "R40" to the literal text "1".
#
# The leading apostrophe is Excel's standard "treat this as text" input
# marker - it forces the numeric-looking value to be stored as a text
# (shared-string) cell instead of being parsed as the number 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "'1"
